$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (names + totals), already sorted descending by total_registros,
# replacing rows 2-12 (rows 13-14 remain unchanged).
$data = @(
    @("PEREZ VEGA ANA YSABEL", 140),
    @("GARAVITO LEON IVONNE LISSETH", 126),
    @("TIMOTEO BAYONA SHARYN LISSETH", 125),
    @("ZAPATA ZETA ROSA ARACELI", 124),
    @("PANTA MONZON SHIRLEY MARIBEL", 111),
    @("NIÑO GUERRERO ANYELA MELINA", 101),
    @("CASTRO JUAREZ MARIA ISABEL", 99),
    @("VALLE SILVA SUTMMER ORFELINDA", 94),
    @("TIZON NUÑEZ FRESIA YAMILI", 90),
    @("MORENO PALACIOS DAMARIS VANESA", 84),
    @("CHERO JUAREZ ANYELA TATIANA", 71)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
